# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Kraken_Profits profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7003.643
$ws.Range("I32").Value = 4263
$ws.Range("K32").Value = 4263
$ws.Range("M32").Value = -3937
$ws.Range("H40").Value = 6439.857
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 6439.857
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 6439.857
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -6789.857
$ws.Range("H64").Value = 3800
$ws.Range("J64").Value = 3800
$ws.Range("L64").Value = 3800
$ws.Range("N64").Value = -4296
$ws.Range("H67").Value = 3800
$ws.Range("J67").Value = 3800
$ws.Range("L67").Value = 3800
$ws.Range("N67").Value = -5516
$ws.Range("H74").Value = 3640.4
$ws.Range("I74").Value = 3600.5
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 3600.5
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -2664.5
$ws.Range("N74").Value = -5672
$ws.Range("H76").Value = 4999.5
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684
$ws.Range("H77").Value = 3640.4
$ws.Range("I77").Value = 3600.5
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 18002.5
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -13322.5
$ws.Range("N77").Value = -28360
$ws.Range("H79").Value = 4999.5
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1440.25
$ws.Range("I61").Value = 1440.25
$ws.Range("K61").Value = 1440.25
$ws.Range("M61").Value = -1228.25
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1094
$ws.Range("N88").ClearContents()
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -96
$ws.Range("N91").ClearContents()
$ws.Range("H132").Value = 2415.2222
$ws.Range("I132").Value = 2306.6667
$ws.Range("K132").Value = 6920.000100000001
$ws.Range("M132").Value = -4390.000100000001
$ws.Range("H136").Value = 1440.25
$ws.Range("I136").Value = 1440.25
$ws.Range("K136").Value = 4320.75
$ws.Range("M136").Value = -1770.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2499.75
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H134").Value = 8000.1113
$ws.Range("I134").Value = 4167.3335
$ws.Range("J134").Value = 15665.667
$ws.Range("K134").Value = 12502.0005
$ws.Range("L134").Value = 46997.001
$ws.Range("M134").Value = -9967.000499999998
$ws.Range("N134").Value = -52067.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 575.5
$ws.Range("I11").Value = 650
$ws.Range("J11").Value = 501
$ws.Range("K11").Value = 650
$ws.Range("L11").Value = 501
$ws.Range("M11").Value = -510
$ws.Range("N11").Value = -781
$ws.Range("H58").Value = 3493.75
$ws.Range("I58").Value = 2991.6667
$ws.Range("K58").Value = 2991.6667
$ws.Range("M58").Value = -2788.6667
$ws.Range("H62").Value = 2966.6667
$ws.Range("I62").Value = 2900
$ws.Range("K62").Value = 2900
$ws.Range("M62").Value = -2276
$ws.Range("H65").Value = 2966.6667
$ws.Range("I65").Value = 2900
$ws.Range("K65").Value = 14500
$ws.Range("M65").Value = -11380
$ws.Range("H107").Value = 282.72726
$ws.Range("I107").Value = 284.66666
$ws.Range("K107").Value = 284.66666
$ws.Range("M107").Value = 1635.33334
$ws.Range("H132").Value = 2254.0833
$ws.Range("I132").Value = 2592.7144
$ws.Range("K132").Value = 7778.1432
$ws.Range("M132").Value = -5248.1432
$ws.Range("H136").Value = 3493.75
$ws.Range("I136").Value = 2991.6667
$ws.Range("K136").Value = 8975.000100000001
$ws.Range("M136").Value = -6425.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2666.6667
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -8888
$ws.Range("N5").Value = -7724
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H134").Value = 9499.5
$ws.Range("I134").Value = 9499.5
$ws.Range("K134").Value = 28498.5
$ws.Range("M134").Value = -23428.5
$ws.Range("H135").Value = 2666.6667
$ws.Range("I135").Value = 3000
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 27000
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -24465
$ws.Range("N135").Value = -27570
$ws.Range("H139").Value = 1913.1666
$ws.Range("I139").Value = 1619.75
$ws.Range("J139").Value = 2500
$ws.Range("K139").Value = 4859.25
$ws.Range("L139").Value = 7500
$ws.Range("M139").Value = 280.75
$ws.Range("N139").Value = -17780
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 4000
$ws.Range("M70").Value = -3730
$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 4000
$ws.Range("M73").Value = -3064
$ws.Range("H80").Value = 2944
$ws.Range("H83").Value = 2944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6500
$ws.Range("I7").Value = 5500
$ws.Range("K7").Value = 5500
$ws.Range("M7").Value = -5388
$ws.Range("H40").Value = 4730.1
$ws.Range("I40").Value = 4287.625
$ws.Range("K40").Value = 4287.625
$ws.Range("M40").Value = -4151.625
$ws.Range("H82").Value = 1784.7
$ws.Range("I82").Value = 1843.375
$ws.Range("K82").Value = 1843.375
$ws.Range("M82").Value = -1482.375
$ws.Range("H85").Value = 1784.7
$ws.Range("I85").Value = 1843.375
$ws.Range("K85").Value = 1843.375
$ws.Range("M85").Value = -595.375
$ws.Range("H122").Value = 5650
$ws.Range("I122").Value = 5499.2856
$ws.Range("J122").Value = 6001.6665
$ws.Range("K122").Value = 16497.8568
$ws.Range("L122").Value = 18004.9995
$ws.Range("M122").Value = -14047.8568
$ws.Range("N122").Value = -22904.9995
$ws.Range("H126").Value = 6500
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030
$ws.Range("H132").Value = 26549.666
$ws.Range("I132").Value = 28657.834
$ws.Range("K132").Value = 85973.50199999999
$ws.Range("M132").Value = -83443.50199999999
$ws.Range("H136").Value = 7665.8335
$ws.Range("I136").Value = 4443.3335
$ws.Range("J136").Value = 17333.334
$ws.Range("K136").Value = 13330.0005
$ws.Range("L136").Value = 52000.00199999999
$ws.Range("M136").Value = -10780.0005
$ws.Range("N136").Value = -57100.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 4023824.5
$ws.Range("I52").Value = 5017530.5
$ws.Range("J52").Value = 49000
$ws.Range("K52").Value = 5017530.5
$ws.Range("L52").Value = 49000
$ws.Range("M52").Value = -5017304.5
$ws.Range("N52").Value = -49452
$ws.Range("H122").Value = 225809.89
$ws.Range("I122").Value = 289042.72
$ws.Range("J122").Value = 4495
$ws.Range("K122").Value = 867128.1599999999
$ws.Range("L122").Value = 13485
$ws.Range("M122").Value = -864678.1599999999
$ws.Range("N122").Value = -18385
$ws.Range("H132").Value = 9359.134
$ws.Range("I132").Value = 5341.143
$ws.Range("J132").Value = 12874.875
$ws.Range("K132").Value = 16023.429
$ws.Range("L132").Value = 38624.625
$ws.Range("M132").Value = -13493.429
$ws.Range("N132").Value = -43684.625
$ws.Range("H136").Value = 1880.75
$ws.Range("I136").Value = 1880.75
$ws.Range("K136").Value = 5642.25
$ws.Range("M136").Value = -3092.25

